$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header cell B1 from "REGION" to "Region"
$ws.Range("B1").Value() = "Region"

# Remove the first three data rows (Negros Occidental / CAUAYAN entries, LMS 2021)
# This shifts all subsequent rows up by 3.
$ws.Rows("14:16").Delete()

# Remove what are now the last three rows (originally rows 38-40, the
# additional Negros Occidental LMS 2024 entries) which are dropped entirely.
$ws.Rows("35:37").Delete()

# Adjust column widths for Division, School Name and Municipality columns
# (ColumnWidth is reported/set ~0.83 below the stored OOXML "width" attribute)
$ws.Columns("C").ColumnWidth = 16.17
$ws.Columns("E").ColumnWidth = 38.17
$ws.Columns("F").ColumnWidth = 19.17
